# Apply "Taking Latest to local from Search Module" edit:
# Adds 9 new test-case rows (101-109) to the "Test Cases" sheet,
# normalises the style of D100, adds a new font/style for B107,
# splits column A/B widths, and updates the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Fix D100 style: s=7 -> s=3 (drop the redundant applyFill flag) ---
$refS3 = $ws.Cells.Item(36, 4)   # D36 already uses the target style (index 3)
$refS3.Copy()
$ws.Cells.Item(100, 4).PasteSpecial(-4122)

# --- Reference cells used to replicate existing cell styles onto the new rows ---
$refS3A = $ws.Cells.Item(36, 1)  # A36 -> style 3 (bordered)
$refS3B = $ws.Cells.Item(36, 2)  # B36 -> style 3 (bordered)
$refS3C = $ws.Cells.Item(36, 3)  # C36 -> style 3 (bordered)
$refS3D = $ws.Cells.Item(36, 4)  # D36 -> style 3 (bordered)
$refS3E = $ws.Cells.Item(36, 5)  # E36 -> style 3 (bordered)
$refS4C = $ws.Cells.Item(71, 3)  # C71 -> style 4 (bordered + wrap text)
$refS6B = $ws.Cells.Item(71, 2)  # B71 -> style 6 (bordered + wrap text + fill-apply)
$refS7A = $ws.Cells.Item(71, 1)  # A71 -> style 7 (bordered + fill-apply)

# --- Cell values, written in the same order the original edit introduced them ---
$ws.Cells.Item(101, 1).Value = "TestCase_B100"
$ws.Cells.Item(101, 2).Value = "OPQA-582"
$ws.Cells.Item(101, 3).Value = "Verify that more search results get displayed when user scrolls down in PATENTS search results page"
$ws.Cells.Item(101, 4).Value = "Y"
$ws.Cells.Item(101, 5).Value = "SKIP"
$ws.Cells.Item(102, 1).Value = "TestCase_B101"
$ws.Cells.Item(102, 2).Value = "OPQA-584"
$ws.Cells.Item(102, 3).Value = "Verify that sorting is retained when user navigates back to PATENTS search results page from record view page"
$ws.Cells.Item(102, 4).Value = "Y"
$ws.Cells.Item(102, 5).Value = "SKIP"
$ws.Cells.Item(103, 1).Value = "TestCase_B102"
$ws.Cells.Item(103, 2).Value = "OPQA-586"
$ws.Cells.Item(103, 3).Value = "Verify that search drop down content type is retained when user navigates back to PATENTS search results page from record view page"
$ws.Cells.Item(103, 4).Value = "Y"
$ws.Cells.Item(103, 5).Value = "SKIP"
$ws.Cells.Item(104, 1).Value = "TestCase_B103"
$ws.Cells.Item(104, 2).Value = "OPQA-591"
$ws.Cells.Item(104, 3).Value = "Verify that filtering is retained when user navigates back to PATENTS search results page from record view page"
$ws.Cells.Item(104, 4).Value = "Y"
$ws.Cells.Item(104, 5).Value = "SKIP"
$ws.Cells.Item(105, 1).Value = "TestCase_B104"
$ws.Cells.Item(106, 1).Value = "TestCase_B105"
$ws.Cells.Item(105, 2).Value = "OPQA-554"
$ws.Cells.Item(105, 3).Value = "Verify that record view page of a post gets displayed when user clicks on article title in ALL  search results page"
$ws.Cells.Item(106, 3).Value = "Verify that record view page of a post gets displayed when user clicks on article title in POSTs search results page`nVerify that following fields get displayed correctly for a post in record view page:`na)Title`n b)Creation date and time `nc)Last edited date and time `nd)Author `ne)Author details `nf)Post content `ng)Likes count `nh)Comments count `ni)Views count"
$ws.Cells.Item(106, 2).Value = "OPQA-555|OPQA-556"
$ws.Cells.Item(105, 4).Value = "Y"
$ws.Cells.Item(105, 5).Value = "SKIP"
$ws.Cells.Item(106, 4).Value = "Y"
$ws.Cells.Item(106, 5).Value = "SKIP"
$ws.Cells.Item(107, 1).Value = "TestCase_B106"
$ws.Cells.Item(107, 3).Value = "Verify that following options get displayed in SORT BY drop down in POSTS search results page: `na)Relevance `nb)Create Date(Newest) `nc)Create Date(Oldest)"
$ws.Cells.Item(107, 2).Value = "OPQA-1226"
$ws.Cells.Item(107, 4).Value = "Y"
$ws.Cells.Item(107, 5).Value = "PASS"
$ws.Cells.Item(108, 1).Value = "TestCase_B107"
$ws.Cells.Item(108, 2).Value = "OPQA-574"
$ws.Cells.Item(108, 3).Value = "Verify that left navigation pane content type is retained when user navigates back to ALL search results page from record view page"
$ws.Cells.Item(108, 4).Value = "Y"
$ws.Cells.Item(108, 5).Value = "SKIP"
$ws.Cells.Item(109, 1).Value = "TestCase_B108"
$ws.Cells.Item(109, 2).Value = "OPQA-569"
$ws.Cells.Item(109, 3).Value = "Verify that sorting is retained when user navigates back to ALL search results page from record view page"
$ws.Cells.Item(109, 4).Value = "Y"
$ws.Cells.Item(109, 5).Value = "PASS"

# --- Cell styles (format-only paste from matching reference cells) ---
# Row 101
$refS3A.Copy()
$ws.Cells.Item(101, 1).PasteSpecial(-4122)
$refS3B.Copy()
$ws.Cells.Item(101, 2).PasteSpecial(-4122)
$refS3C.Copy()
$ws.Cells.Item(101, 3).PasteSpecial(-4122)
$refS3D.Copy()
$ws.Cells.Item(101, 4).PasteSpecial(-4122)
$refS3E.Copy()
$ws.Cells.Item(101, 5).PasteSpecial(-4122)
# Row 102
$refS3A.Copy()
$ws.Cells.Item(102, 1).PasteSpecial(-4122)
$refS3B.Copy()
$ws.Cells.Item(102, 2).PasteSpecial(-4122)
$refS3C.Copy()
$ws.Cells.Item(102, 3).PasteSpecial(-4122)
$refS3D.Copy()
$ws.Cells.Item(102, 4).PasteSpecial(-4122)
$refS3E.Copy()
$ws.Cells.Item(102, 5).PasteSpecial(-4122)
# Row 103
$refS3A.Copy()
$ws.Cells.Item(103, 1).PasteSpecial(-4122)
$refS3B.Copy()
$ws.Cells.Item(103, 2).PasteSpecial(-4122)
$refS3C.Copy()
$ws.Cells.Item(103, 3).PasteSpecial(-4122)
$refS3D.Copy()
$ws.Cells.Item(103, 4).PasteSpecial(-4122)
$refS3E.Copy()
$ws.Cells.Item(103, 5).PasteSpecial(-4122)
# Row 104
$refS3A.Copy()
$ws.Cells.Item(104, 1).PasteSpecial(-4122)
$refS3B.Copy()
$ws.Cells.Item(104, 2).PasteSpecial(-4122)
$refS3C.Copy()
$ws.Cells.Item(104, 3).PasteSpecial(-4122)
$refS3D.Copy()
$ws.Cells.Item(104, 4).PasteSpecial(-4122)
$refS3E.Copy()
$ws.Cells.Item(104, 5).PasteSpecial(-4122)
# Row 105
$refS3A.Copy()
$ws.Cells.Item(105, 1).PasteSpecial(-4122)
$refS3B.Copy()
$ws.Cells.Item(105, 2).PasteSpecial(-4122)
$refS3C.Copy()
$ws.Cells.Item(105, 3).PasteSpecial(-4122)
$refS3D.Copy()
$ws.Cells.Item(105, 4).PasteSpecial(-4122)
$refS3E.Copy()
$ws.Cells.Item(105, 5).PasteSpecial(-4122)
# Row 106
$refS3A.Copy()
$ws.Cells.Item(106, 1).PasteSpecial(-4122)
$refS3B.Copy()
$ws.Cells.Item(106, 2).PasteSpecial(-4122)
$refS4C.Copy()
$ws.Cells.Item(106, 3).PasteSpecial(-4122)
$refS3D.Copy()
$ws.Cells.Item(106, 4).PasteSpecial(-4122)
$refS3E.Copy()
$ws.Cells.Item(106, 5).PasteSpecial(-4122)
# Row 107
$refS3A.Copy()
$ws.Cells.Item(107, 1).PasteSpecial(-4122)
$refS7A.Copy()
$ws.Cells.Item(107, 2).PasteSpecial(-4122)
$ws.Cells.Item(107, 2).Font.Color = 0
$refS4C.Copy()
$ws.Cells.Item(107, 3).PasteSpecial(-4122)
$refS3D.Copy()
$ws.Cells.Item(107, 4).PasteSpecial(-4122)
$refS3E.Copy()
$ws.Cells.Item(107, 5).PasteSpecial(-4122)
# Row 108
$refS7A.Copy()
$ws.Cells.Item(108, 1).PasteSpecial(-4122)
$refS6B.Copy()
$ws.Cells.Item(108, 2).PasteSpecial(-4122)
$refS4C.Copy()
$ws.Cells.Item(108, 3).PasteSpecial(-4122)
$refS3D.Copy()
$ws.Cells.Item(108, 4).PasteSpecial(-4122)
$refS3E.Copy()
$ws.Cells.Item(108, 5).PasteSpecial(-4122)
# Row 109
$refS7A.Copy()
$ws.Cells.Item(109, 1).PasteSpecial(-4122)
$refS6B.Copy()
$ws.Cells.Item(109, 2).PasteSpecial(-4122)
$refS4C.Copy()
$ws.Cells.Item(109, 3).PasteSpecial(-4122)
$refS3D.Copy()
$ws.Cells.Item(109, 4).PasteSpecial(-4122)
$refS3E.Copy()
$ws.Cells.Item(109, 5).PasteSpecial(-4122)

# --- Explicit row heights (11 wrapped lines x15 / 4 wrapped lines x15) ---
$ws.Rows.Item(106).RowHeight = 165
$ws.Rows.Item(107).RowHeight = 60

# --- Column widths: column B now has its own (wider) best-fit width ---
$ws.Columns.Item(2).AutoFit()

$excel.CutCopyMode = $false

# --- Restore the final on-screen selection ---
$ws.Range("C122").Select()

